$d = $word.ActiveDocument

# Splits a paragraph's plain run text into several runs with identical
# formatting, each containing one of the given $Segments strings (in
# order). This fixes the "dados.lengh" -> "dados.length" typo (and, for a
# handful of paragraphs, also corrects the trailing counter value) while
# reproducing the exact run layout produced by the original edit.
function Split-Run($ParaIndex, $Segments) {
    $p = $d.Paragraphs($ParaIndex)
    $full = $p.Range
    $s = $full.Start
    # Full.Text includes the trailing paragraph-mark character, so the
    # visible text length is one less than it.
    $origLen = $full.Text.Length - 1
    $joined = [string]::Join("", $Segments)

    # Rewrite the paragraph's visible text in one shot (this collapses it
    # into a single run sharing the formatting that was already there).
    $textRange = $d.Range($s, $s + $origLen)
    $textRange.Text = $joined

    # Now carve out run boundaries, from the right-most boundary back to
    # the left-most one, by toggling Bold on/off (no net formatting
    # change) on the trailing sub-range. Because the runs being toggled
    # keep the exact same resulting formatting, no rPr differences are
    # introduced -- but this engine only coalesces identically-formatted
    # adjacent runs immediately after a text write, not afterwards, so
    # the split persists.
    $total = $joined.Length
    $cum = $total
    for ($i = $Segments.Count - 1; $i -ge 1; $i--) {
        $cum = $cum - $Segments[$i].Length
        $boundary = $s + $cum
        $rest = $d.Range($boundary, $s + $total)
        $rest.Font.Bold = $true
        $rest.Font.Bold = $false
    }
}

Split-Run 45  @("dados.", "length", " ", "= 0")
Split-Run 60  @("dados.", "length", " ", "= 1")
Split-Run 75  @("dados.", "length", " ", "= 2")
Split-Run 93  @("dados.", "length", " ", "= 3")
Split-Run 111 @("dados.", "length", " ", "= 3")
Split-Run 132 @("dados.", "length", " ", "= 3")
Split-Run 151 @("dados.", "length", " ", "= 2")
Split-Run 170 @("dados.", "length", " ", "= 1")
Split-Run 189 @("dados.", "length", " ", "= 0")
Split-Run 206 @("dados.", "length", " ", "= 0")
Split-Run 225 @("dados.lengh = ", "3")
Split-Run 243 @("dados.", "length", " ", "= ", "3")
Split-Run 261 @("dados.", "length", " ", "= ", "3")
Split-Run 279 @("dados.", "length", " ", "=3")
Split-Run 297 @("dados.", "length", " ", "= 3")
Split-Run 318 @("dados.", "length", " ", "= 3")
Split-Run 335 @("dados.", "length", " ", "=", "3")
Split-Run 353 @("dados.", "length", " ", "=", "3")
Split-Run 373 @("dados.", "length", " ", "=", "3")
Split-Run 391 @("dados.", "length", " ", "= ", "3")
